$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("model")
$ws2.Activate()
$ws2.Range("B5").Select()

$ws = $wb.Worksheets.Item("survey")
$ws.Rows.Item(29).Resize(2).Delete()
$ws.Select()
$ws.Range("A29:XFD30").Select()
